# Refactor synthetic array /3 (publipostage fix)
# The shared strings used as "status icon" / "status name" values are being
# renamed:
#   ⬛ -> 📘   (noir -> bleu)
#   🟩 -> 📗
#   🟧 -> 📙
#   🟥 -> 📕
#   noir -> bleu
#
# We apply this by rewriting every cell in columns A (icon) and B (label)
# that currently holds one of the old values with its new counterpart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$iconMap = @{
    "⬛" = "📘"
    "🟩" = "📗"
    "🟧" = "📙"
    "🟥" = "📕"
}

$labelMap = @{
    "noir" = "bleu"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value()
    if ($iconMap.ContainsKey($aVal)) {
        $aCell.Value = $iconMap[$aVal]
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value()
    if ($labelMap.ContainsKey($bVal)) {
        $bCell.Value = $labelMap[$bVal]
    }
}
